# Swap the species-record data between row 4 and row 5.
# (The underlying DB ids changed, so the two records' content
#  effectively traded places.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $rng4 = $ws.Range($col + "4")
    $rng5 = $ws.Range($col + "5")
    $v4 = $rng4.Value2
    $v5 = $rng5.Value2
    $rng4.Value = $v5
    $rng5.Value = $v4
}

# Column AC ("Publik kommentar") currently only has a value on row 5.
# After the edit it belongs to row 4 and row 5's cell becomes empty.
$ws.Range("AC4").Value = "Senvuxen"
$ws.Range("AC5").Value = $null
